$wb = $excel.ActiveWorkbook

# "Repayment Schedule" sheet: a new column is inserted before column N
# (shifting the old "Late" column N -> O and "Outstanding" column P -> Q).
$ws = $wb.Worksheets.Item("Repayment Schedule")
$ws.Columns("N:N").Insert()

# Selection on the Repayment Schedule sheet moves to R7, and that sheet
# becomes the active tab (previously "Transactions" was active/selected).
[void]$ws.Range("R7").Select()
$ws.Activate()
